$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-02-04 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-02-05 Monday", 2) | Out-Null
$d.Content.Find.Execute("26+66=92", $true, $false, $false, $false, $false, $true, 1, $false, "1+12=13", 2) | Out-Null
$d.Content.Find.Execute("37-31=6", $true, $false, $false, $false, $false, $true, 1, $false, "55+32=87", 2) | Out-Null
$d.Content.Find.Execute("35-7=28", $true, $false, $false, $false, $false, $true, 1, $false, "76-73=3", 2) | Out-Null
$d.Content.Find.Execute("36-20=16", $true, $false, $false, $false, $false, $true, 1, $false, "13+78=91", 2) | Out-Null
$d.Content.Find.Execute("4+24=28", $true, $false, $false, $false, $false, $true, 1, $false, "48+11=59", 2) | Out-Null
$d.Content.Find.Execute("86-67=19", $true, $false, $false, $false, $false, $true, 1, $false, "61-1=60", 2) | Out-Null
$d.Content.Find.Execute("24+60=84", $true, $false, $false, $false, $false, $true, 1, $false, "89-40=49", 2) | Out-Null
$d.Content.Find.Execute("0+15=15", $true, $false, $false, $false, $false, $true, 1, $false, "32+47=79", 2) | Out-Null
$d.Content.Find.Execute("90-20=70", $true, $false, $false, $false, $false, $true, 1, $false, "32-17=15", 2) | Out-Null
$d.Content.Find.Execute("82-65=17", $true, $false, $false, $false, $false, $true, 1, $false, "98-32=66", 2) | Out-Null
$d.Content.Find.Execute("69-27=42", $true, $false, $false, $false, $false, $true, 1, $false, "43+9=52", 2) | Out-Null
$d.Content.Find.Execute("43+51=94", $true, $false, $false, $false, $false, $true, 1, $false, "4+76=80", 2) | Out-Null
$d.Content.Find.Execute("68-0=68", $true, $false, $false, $false, $false, $true, 1, $false, "6-2=4", 2) | Out-Null
$d.Content.Find.Execute("68-8=60", $true, $false, $false, $false, $false, $true, 1, $false, "65-23=42", 2) | Out-Null
$d.Content.Find.Execute("77-18=59", $true, $false, $false, $false, $false, $true, 1, $false, "1+7=8", 2) | Out-Null
$d.Content.Find.Execute("81+0=81", $true, $false, $false, $false, $false, $true, 1, $false, "12+29=41", 2) | Out-Null
$d.Content.Find.Execute("1+47=48", $true, $false, $false, $false, $false, $true, 1, $false, "43+45=88", 2) | Out-Null
$d.Content.Find.Execute("70+18=88", $true, $false, $false, $false, $false, $true, 1, $false, "95-26=69", 2) | Out-Null
$d.Content.Find.Execute("95-36=59", $true, $false, $false, $false, $false, $true, 1, $false, "69-44=25", 2) | Out-Null
$d.Content.Find.Execute("23+46=69", $true, $false, $false, $false, $false, $true, 1, $false, "38+52=90", 2) | Out-Null
$d.Content.Find.Execute("73-25=48", $true, $false, $false, $false, $false, $true, 1, $false, "86-78=8", 2) | Out-Null
$d.Content.Find.Execute("55+24=79", $true, $false, $false, $false, $false, $true, 1, $false, "43+49=92", 2) | Out-Null
$d.Content.Find.Execute("90-35=55", $true, $false, $false, $false, $false, $true, 1, $false, "64+6=70", 2) | Out-Null
$d.Content.Find.Execute("95-45=50", $true, $false, $false, $false, $false, $true, 1, $false, "64+14=78", 2) | Out-Null
$d.Content.Find.Execute("49-37=12", $true, $false, $false, $false, $false, $true, 1, $false, "28+35=63", 2) | Out-Null
$d.Content.Find.Execute("57+38=95", $true, $false, $false, $false, $false, $true, 1, $false, "52+0=52", 2) | Out-Null
$d.Content.Find.Execute("87+8=95", $true, $false, $false, $false, $false, $true, 1, $false, "23+20=43", 2) | Out-Null
$d.Content.Find.Execute("49+15=64", $true, $false, $false, $false, $false, $true, 1, $false, "24-22=2", 2) | Out-Null
$d.Content.Find.Execute("60+33=93", $true, $false, $false, $false, $false, $true, 1, $false, "44-34=10", 2) | Out-Null
$d.Content.Find.Execute("23+42=65", $true, $false, $false, $false, $false, $true, 1, $false, "85-77=8", 2) | Out-Null
$d.Content.Find.Execute("55-8=47", $true, $false, $false, $false, $false, $true, 1, $false, "98-29=69", 2) | Out-Null
$d.Content.Find.Execute("83-73=10", $true, $false, $false, $false, $false, $true, 1, $false, "32-3=29", 2) | Out-Null
$d.Content.Find.Execute("9+73=82", $true, $false, $false, $false, $false, $true, 1, $false, "83-71=12", 2) | Out-Null
$d.Content.Find.Execute("17+57=74", $true, $false, $false, $false, $false, $true, 1, $false, "73-61=12", 2) | Out-Null
$d.Content.Find.Execute("38-8=30", $true, $false, $false, $false, $false, $true, 1, $false, "35+35=70", 2) | Out-Null
$d.Content.Find.Execute("55-16=39", $true, $false, $false, $false, $false, $true, 1, $false, "56+6=62", 2) | Out-Null
$d.Content.Find.Execute("51+36=87", $true, $false, $false, $false, $false, $true, 1, $false, "92-11=81", 2) | Out-Null
$d.Content.Find.Execute("3+11=14", $true, $false, $false, $false, $false, $true, 1, $false, "81-17=64", 2) | Out-Null
$d.Content.Find.Execute("17-9=8", $true, $false, $false, $false, $false, $true, 1, $false, "45-38=7", 2) | Out-Null
$d.Content.Find.Execute("98-43=55", $true, $false, $false, $false, $false, $true, 1, $false, "39-21=18", 2) | Out-Null
$d.Content.Find.Execute("35+61=96", $true, $false, $false, $false, $false, $true, 1, $false, "35+33=68", 2) | Out-Null
$d.Content.Find.Execute("28+0=28", $true, $false, $false, $false, $false, $true, 1, $false, "51-4=47", 2) | Out-Null
$d.Content.Find.Execute("77-35=42", $true, $false, $false, $false, $false, $true, 1, $false, "43-13=30", 2) | Out-Null
$d.Content.Find.Execute("67-56=11", $true, $false, $false, $false, $false, $true, 1, $false, "77-24=53", 2) | Out-Null
$d.Content.Find.Execute("96-43=53", $true, $false, $false, $false, $false, $true, 1, $false, "12+33=45", 2) | Out-Null
$d.Content.Find.Execute("91-58=33", $true, $false, $false, $false, $false, $true, 1, $false, "68+17=85", 2) | Out-Null
$d.Content.Find.Execute("49+33=82", $true, $false, $false, $false, $false, $true, 1, $false, "85+1=86", 2) | Out-Null
$d.Content.Find.Execute("55+12=67", $true, $false, $false, $false, $false, $true, 1, $false, "32+34=66", 2) | Out-Null
$d.Content.Find.Execute("45-33=12", $true, $false, $false, $false, $false, $true, 1, $false, "71+23=94", 2) | Out-Null
$d.Content.Find.Execute("50-6=44", $true, $false, $false, $false, $false, $true, 1, $false, "60-0=60", 2) | Out-Null
$d.Content.Find.Execute("44+53=97", $true, $false, $false, $false, $false, $true, 1, $false, "58+34=92", 2) | Out-Null
$d.Content.Find.Execute("29+68=97", $true, $false, $false, $false, $false, $true, 1, $false, "46-8=38", 2) | Out-Null
$d.Content.Find.Execute("40+50=90", $true, $false, $false, $false, $false, $true, 1, $false, "31+33=64", 2) | Out-Null
$d.Content.Find.Execute("87-85=2", $true, $false, $false, $false, $false, $true, 1, $false, "80+13=93", 2) | Out-Null
$d.Content.Find.Execute("7+8=15", $true, $false, $false, $false, $false, $true, 1, $false, "81-49=32", 2) | Out-Null
$d.Content.Find.Execute("74-2=72", $true, $false, $false, $false, $false, $true, 1, $false, "37+35=72", 2) | Out-Null
$d.Content.Find.Execute("58-10=48", $true, $false, $false, $false, $false, $true, 1, $false, "0+55=55", 2) | Out-Null
$d.Content.Find.Execute("29+41=70", $true, $false, $false, $false, $false, $true, 1, $false, "72-48=24", 2) | Out-Null
$d.Content.Find.Execute("36-6=30", $true, $false, $false, $false, $false, $true, 1, $false, "34+47=81", 2) | Out-Null
$d.Content.Find.Execute("10+72=82", $true, $false, $false, $false, $false, $true, 1, $false, "9+54=63", 2) | Out-Null
$d.Content.Find.Execute("19-14=5", $true, $false, $false, $false, $false, $true, 1, $false, "37-0=37", 2) | Out-Null
$d.Content.Find.Execute("90+1=91", $true, $false, $false, $false, $false, $true, 1, $false, "7+64=71", 2) | Out-Null
$d.Content.Find.Execute("38-35=3", $true, $false, $false, $false, $false, $true, 1, $false, "90-43=47", 2) | Out-Null
$d.Content.Find.Execute("90-84=6", $true, $false, $false, $false, $false, $true, 1, $false, "37+3=40", 2) | Out-Null
$d.Content.Find.Execute("61+2=63", $true, $false, $false, $false, $false, $true, 1, $false, "63+6=69", 2) | Out-Null
$d.Content.Find.Execute("69-16=53", $true, $false, $false, $false, $false, $true, 1, $false, "59-49=10", 2) | Out-Null
$d.Content.Find.Execute("99-50=49", $true, $false, $false, $false, $false, $true, 1, $false, "50+26=76", 2) | Out-Null
$d.Content.Find.Execute("94-73=21", $true, $false, $false, $false, $false, $true, 1, $false, "2+22=24", 2) | Out-Null
$d.Content.Find.Execute("90-37=53", $true, $false, $false, $false, $false, $true, 1, $false, "89-13=76", 2) | Out-Null
$d.Content.Find.Execute("6+86=92", $true, $false, $false, $false, $false, $true, 1, $false, "18+7=25", 2) | Out-Null
$d.Content.Find.Execute("84-24=60", $true, $false, $false, $false, $false, $true, 1, $false, "50-13=37", 2) | Out-Null
$d.Content.Find.Execute("33-10=23", $true, $false, $false, $false, $false, $true, 1, $false, "69+19=88", 2) | Out-Null
$d.Content.Find.Execute("49-13=36", $true, $false, $false, $false, $false, $true, 1, $false, "91+0=91", 2) | Out-Null
$d.Content.Find.Execute("45-26=19", $true, $false, $false, $false, $false, $true, 1, $false, "36+23=59", 2) | Out-Null
$d.Content.Find.Execute("83-30=53", $true, $false, $false, $false, $false, $true, 1, $false, "26-8=18", 2) | Out-Null
$d.Content.Find.Execute("61-22=39", $true, $false, $false, $false, $false, $true, 1, $false, "26+52=78", 2) | Out-Null
$d.Content.Find.Execute("38-11=27", $true, $false, $false, $false, $false, $true, 1, $false, "64-53=11", 2) | Out-Null
$d.Content.Find.Execute("17+4=21", $true, $false, $false, $false, $false, $true, 1, $false, "43-20=23", 2) | Out-Null
$d.Content.Find.Execute("64-43=21", $true, $false, $false, $false, $false, $true, 1, $false, "13+8=21", 2) | Out-Null
$d.Content.Find.Execute("71+8=79", $true, $false, $false, $false, $false, $true, 1, $false, "54+15=69", 2) | Out-Null
$d.Content.Find.Execute("40+27=67", $true, $false, $false, $false, $false, $true, 1, $false, "21+67=88", 2) | Out-Null
$d.Content.Find.Execute("60-32=28", $true, $false, $false, $false, $false, $true, 1, $false, "96-50=46", 2) | Out-Null
$d.Content.Find.Execute("92-90=2", $true, $false, $false, $false, $false, $true, 1, $false, "19+31=50", 2) | Out-Null
$d.Content.Find.Execute("97-74=23", $true, $false, $false, $false, $false, $true, 1, $false, "34+41=75", 2) | Out-Null
$d.Content.Find.Execute("35+54=89", $true, $false, $false, $false, $false, $true, 1, $false, "6+62=68", 2) | Out-Null
$d.Content.Find.Execute("73-56=17", $true, $false, $false, $false, $false, $true, 1, $false, "24-9=15", 2) | Out-Null
$d.Content.Find.Execute("37-25=12", $true, $false, $false, $false, $false, $true, 1, $false, "22+48=70", 2) | Out-Null
$d.Content.Find.Execute("67-45=22", $true, $false, $false, $false, $false, $true, 1, $false, "57+34=91", 2) | Out-Null
$d.Content.Find.Execute("29+6=35", $true, $false, $false, $false, $false, $true, 1, $false, "0+32=32", 2) | Out-Null
$d.Content.Find.Execute("56-53=3", $true, $false, $false, $false, $false, $true, 1, $false, "94-93=1", 2) | Out-Null
$d.Content.Find.Execute("99-41=58", $true, $false, $false, $false, $false, $true, 1, $false, "52-31=21", 2) | Out-Null
$d.Content.Find.Execute("21+26=47", $true, $false, $false, $false, $false, $true, 1, $false, "8+32=40", 2) | Out-Null
$d.Content.Find.Execute("62-49=13", $true, $false, $false, $false, $false, $true, 1, $false, "8+54=62", 2) | Out-Null
$d.Content.Find.Execute("86-73=13", $true, $false, $false, $false, $false, $true, 1, $false, "17+64=81", 2) | Out-Null
$d.Content.Find.Execute("55-25=30", $true, $false, $false, $false, $false, $true, 1, $false, "37+49=86", 2) | Out-Null
$d.Content.Find.Execute("16+46=62", $true, $false, $false, $false, $false, $true, 1, $false, "44+13=57", 2) | Out-Null
$d.Content.Find.Execute("90+7=97", $true, $false, $false, $false, $false, $true, 1, $false, "29+1=30", 2) | Out-Null
$d.Content.Find.Execute("19-1=18", $true, $false, $false, $false, $false, $true, 1, $false, "64-34=30", 2) | Out-Null
$d.Content.Find.Execute("45+46=91", $true, $false, $false, $false, $false, $true, 1, $false, "7+2=9", 2) | Out-Null
$d.Content.Find.Execute("68-10=58", $true, $false, $false, $false, $false, $true, 1, $false, "18+61=79", 2) | Out-Null

Write-Output "Replacements applied: 101"
